$wb = $excel.ActiveWorkbook

# The same F-column updates need to be applied to both the "展览" sheet
# and the "全部类型" sheet, which mirror the same underlying data.
$sheetNames = @("展览", "全部类型")

# Map of row -> new value for column F ("想去人数")
$updates = @{
    6  = 761
    9  = 4511
    11 = 358
    12 = 1287
    13 = 536
    15 = 875
    17 = 483
    19 = 231
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
